# Updates benchmark-result values in the single-column results table.
# The table has one row per measurement; Cell.Range.Text replaces only the
# text content of the cell while leaving the existing run formatting
# (Times New Roman, sz 22) untouched, matching the target XML exactly.
# (A couple of rows -- 44/45/46 -- originally hold several tab-separated
# runs summarizing a per-iteration timing breakdown; setting .Range.Text
# collapses them into the single plain-text summary value required.)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-RowValue($rowIndex, $expectedOld, $newValue) {
    $cell = $t.Cell($rowIndex, 1)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $expectedOld) {
        Write-Host "Warning: row $rowIndex expected '$expectedOld' but found '$current'"
    }
    $cell.Range.Text = $newValue
}

Set-RowValue 1  "100"     "0M"
Set-RowValue 2  "0"       "0M"
Set-RowValue 3  "229"     "0M"
Set-RowValue 4  "13"      "55"
Set-RowValue 6  "0.00011" "0.00012"
Set-RowValue 9  "0.00007" "0.00008"
Set-RowValue 10 "0.00007" "0.00009"
Set-RowValue 11 "0.00008" "0.00010"
Set-RowValue 12 "0.00097" "0.00473"
Set-RowValue 44 "2`t0.00007`t0.00007`t0.00007`t0.00000`t0.00007`t0.00007`t0.00007`t0.00014`t100.0" "100"
Set-RowValue 45 "2`t0.00005`t0.00010`t0.00007`t0.00003`t0.00005`t0.00005`t0.00010`t0.00015`t100.0" "0"
Set-RowValue 46 "38`t0.00006`t0.00012`t0.00009`t0.00001`t0.00008`t0.00009`t0.00010`t0.00347`t100.0" "229"
